$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Change 1: remove the "_GoBack" bookmark from the paragraph that reads
# "!place accordingly in folder hierarchy "
# -------------------------------------------------------------------------
$pFolder = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*place accordingly in folder hierarchy*") {
        $pFolder = $p
        break
    }
}

if ($pFolder -ne $null) {
    $xmlFolder = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="070688D5" w14:textId="2272C5F3" w:rsidR="003C55D7" w:rsidRPr="00186863" w:rsidRDefault="003C55D7" w:rsidP="00186863" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:rFonts w:eastAsiaTheme="minorHAnsi" w:cstheme="minorBidi"/><w:lang w:val="en-GB" w:eastAsia="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00186863"><w:rPr><w:rFonts w:eastAsiaTheme="minorHAnsi"/><w:lang w:val="en-GB" w:eastAsia="en-US"/></w:rPr><w:t xml:space="preserve">!place accordingly in folder hierarchy </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$pFolder.Range.InsertXML($xmlFolder)
}

# -------------------------------------------------------------------------
# Change 2: the empty paragraph right after "... ('create_jsonfile_bold'
# or '_T1w')" becomes a new paragraph of text about the echo-spacing excel
# file, followed by a new (still empty) paragraph that now carries the
# "_GoBack" bookmark that used to live on the "folder hierarchy" paragraph.
# -------------------------------------------------------------------------
$pEmpty = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*create_jsonfile_bold*") {
        $pEmpty = $p.Next(1)
        break
    }
}

if ($pEmpty -ne $null) {
    $xmlEmpty = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="30ACE33A" w14:textId="7395EAE8" w:rsidR="007D1FE3" w:rsidRDefault="007D1FE3" w:rsidP="00420DF1" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:r><w:t xml:space="preserve">An excel file called echo spacing contains formulas to calculate effective echo spacing and total readout time, for the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>func</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> sequence runs (you need ETL = EPI factor, WFS, SENSE or whatever factor (MB factor can be ignored); these things you can get from a .txt file that you can create from the sequence parameters on the MR scanner computer) </w:t></w:r></w:p><w:p w14:paraId="30ACE33B" w14:textId="7395EAE9" w:rsidR="007D1FE3" w:rsidRDefault="007D1FE3" w:rsidP="00420DF1" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$pEmpty.Range.InsertXML($xmlEmpty)
}

# -------------------------------------------------------------------------
# Change 3: "!If you have DICOMs available, use dcm2niix ..." becomes
# "!If you have ParRec or Dicom available, use dcm2niix ..."
# -------------------------------------------------------------------------
$pDicom = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*DICOMs available*") {
        $pDicom = $p
        break
    }
}

if ($pDicom -ne $null) {
    $xmlDicom = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="5961AC85" w14:textId="4D1200BC" w:rsidR="007D1FE3" w:rsidRDefault="007D1FE3" w:rsidP="00420DF1" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:proofErr w:type="gramStart"/><w:r><w:t>!If</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> you have </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ParRec</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> or </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Dicom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> available, use </w:t></w:r><w:r w:rsidRPr="007D1FE3"><w:t>dcm2niix</w:t></w:r><w:r><w:t xml:space="preserve"> since this will create both </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nifti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and json files for you</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$pDicom.Range.InsertXML($xmlDicom)
}

Write-Host "Edits applied"
